# "Updated cryptos list on Thu Aug 10 20:06:00 UTC 2023 with GitHub Actions"
#
# Refreshes the Price (column D) and Volume(1h) (column E) columns of the
# crypto snapshot table with the latest coinranking.com figures, and
# re-ranks WrappedEther / Dogecoin / Cardano (rows 8-10) to reflect their
# new relative ordering (Name + Link move together with the row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "29.405.72"; DNumeric = $false; E = "  -0.07%  " },
    @{ Row = 3; D = "1.847.34"; DNumeric = $false; E = "  -0.07%  " },
    @{ Row = 4; D = "0.9987"; DNumeric = $true; E = "  -0.08%  " },
    @{ Row = 5; D = "240.67"; DNumeric = $true; E = "  -0.98%  " },
    @{ Row = 6; D = "0.6333"; DNumeric = $true; E = "  -4.03%  " },
    @{ Row = 7; D = "0.9997"; DNumeric = $true },
    @{ Row = 8; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "3.004.78"; DNumeric = $false; E = "  +62.65%  " },
    @{ Row = 9; B = "Dogecoin"; C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D = "0.07574"; DNumeric = $true; E = "  +1.35%  " },
    @{ Row = 10; B = "Cardano"; C = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D = "0.2967"; DNumeric = $true; E = "  -0.74%  " },
    @{ Row = 11; D = "24.54"; DNumeric = $true; E = "  +0.91%  " },
    @{ Row = 12; D = "0.07712"; DNumeric = $true; E = "  +1.03%  " },
    @{ Row = 13; E = "  -0.62%  " },
    @{ Row = 14; D = "0.6854"; DNumeric = $true; E = "  +0.21%  " },
    @{ Row = 15; D = "82.78"; DNumeric = $true; E = "  -1.03%  " },
    @{ Row = 16; D = "0.000009930"; DNumeric = $true; E = "  +4.60%  " },
    @{ Row = 17; D = "6.184"; DNumeric = $true; E = "  +0.77%  " },
    @{ Row = 18; D = "29.441.91"; DNumeric = $false; E = "  -0.09%  " },
    @{ Row = 19; D = "231.79"; DNumeric = $true; E = "  -2.22%  " },
    @{ Row = 20; D = "12.50"; DNumeric = $true; E = "  -0.47%  " },
    @{ Row = 21; D = "0.9997"; DNumeric = $true; E = "  -0.05%  " },
    @{ Row = 22; D = "7.571"; DNumeric = $true; E = "  -1.30%  " },
    @{ Row = 23; D = "0.9998"; DNumeric = $true },
    @{ Row = 24; D = "154.88"; DNumeric = $true; E = "  -1.24%  " },
    @{ Row = 25; D = "0.1390"; DNumeric = $true; E = "  -2.43%  " },
    @{ Row = 26; D = "8.432"; DNumeric = $true; E = "  -0.54%  " },
    @{ Row = 27; D = "17.66"; DNumeric = $true; E = "  -0.75%  " },
    @{ Row = 28; D = "1.473"; DNumeric = $true; E = "  -1.02%  " },
    @{ Row = 29; D = "0.05803"; DNumeric = $true; E = "  -3.69%  " },
    @{ Row = 30; D = "1.259"; DNumeric = $true; E = "  +0.60%  " },
    @{ Row = 31; E = "  -0.35%  " },
    @{ Row = 32; D = "4.019"; DNumeric = $true; E = "  -1.23%  " },
    @{ Row = 33; D = "3.145.37"; DNumeric = $false; E = "  +57.26%  " },
    @{ Row = 34; E = "  +0.84%  " },
    @{ Row = 35; D = "1.159"; DNumeric = $true; E = "  -1.52%  " },
    @{ Row = 36; D = "0.7189"; DNumeric = $true; E = "  -0.17%  " },
    @{ Row = 37; D = "2.596"; DNumeric = $true; E = "  +0.01%  " },
    @{ Row = 38; D = "1.248.95"; DNumeric = $false; E = "  +4.40%  " },
    @{ Row = 39; D = "2.793"; DNumeric = $true; E = "  -0.17%  " },
    @{ Row = 40; D = "0.01805"; DNumeric = $true; E = "  +1.42%  " },
    @{ Row = 41; D = "0.9057"; DNumeric = $true; E = "  -0.86%  " },
    @{ Row = 42; D = "6.071"; DNumeric = $true; E = "  -2.54%  " },
    @{ Row = 43; D = "0.9989"; DNumeric = $true; E = "  -0.09%  " },
    @{ Row = 44; D = "101.32"; DNumeric = $true },
    @{ Row = 45; D = "67.01"; DNumeric = $true; E = "  +1.33%  " },
    @{ Row = 46; D = "7.321"; DNumeric = $true; E = "  -1.41%  " },
    @{ Row = 47; D = "9.177"; DNumeric = $true; E = "  +1.01%  " },
    @{ Row = 48; E = "  -0.96%  " },
    @{ Row = 49; D = "1.695"; DNumeric = $true; E = "  +2.65%  " },
    @{ Row = 50; D = "0.1124"; DNumeric = $true; E = "  -0.13%  " },
    @{ Row = 51; D = "0.05744"; DNumeric = $true; E = "  +0.03%  " }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($u.ContainsKey("B")) { $ws.Range("B$r").Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C$r").Value = $u.C }

    if ($u.ContainsKey("D")) {
        $dCell = $ws.Range("D$r")
        if ($u.DNumeric) {
            # The feed's Price column is plain text (e.g. "0.9987", "12.50",
            # "29.405.72" grouped-thousands). Values that parse as a plain
            # number (e.g. "0.9987") need the cell forced to Text first,
            # otherwise Excel would silently convert them to a numeric value
            # and drop significant trailing zeros / formatting.
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }

    if ($u.ContainsKey("E")) { $ws.Range("E$r").Value = $u.E }
}
